$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Resolving-Mac rows (sending cluster rows 22-25)
$ws.Rows("22:25").Delete()

# Update recalculated TPM-derived values for rows 2-21
$ws.Cells.Item(2,7).Value = 1.597544666666667
$ws.Cells.Item(2,8).Value = 4.792634
$ws.Cells.Item(2,9).Value = 0.02295281024579055
$ws.Cells.Item(2,10).Value = 0.02295281024579055
$ws.Cells.Item(2,13).Value = 2.680128333333334
$ws.Cells.Item(2,14).Value = 8.040385000000001
$ws.Cells.Item(2,15).Value = 0.8195172133182701
$ws.Cells.Item(2,16).Value = 0.81951721331827
$ws.Cells.Item(2,17).Value = 4.281624724898889
$ws.Cells.Item(2,18).Value = 38.53462252409
$ws.Cells.Item(2,19).Value = 0.01881022309045331
$ws.Cells.Item(2,20).Value = 0.01881022309045331
$ws.Cells.Item(3,7).Value = 1.597544666666667
$ws.Cells.Item(3,8).Value = 4.792634
$ws.Cells.Item(3,9).Value = 0.02295281024579055
$ws.Cells.Item(3,10).Value = 0.02295281024579055
$ws.Cells.Item(3,13).Value = 0.3764886666666666
$ws.Cells.Item(3,15).Value = 0.1151209586179932
$ws.Cells.Item(3,16).Value = 0.1151209586179932
$ws.Cells.Item(3,17).Value = 0.6014574614937778
$ws.Cells.Item(3,18).Value = 5.413117153443999
$ws.Cells.Item(3,19).Value = 0.002642349518472304
$ws.Cells.Item(3,20).Value = 0.002642349518472303
$ws.Cells.Item(4,7).Value = 1.597544666666667
$ws.Cells.Item(4,8).Value = 4.792634
$ws.Cells.Item(4,9).Value = 0.02295281024579055
$ws.Cells.Item(4,10).Value = 0.02295281024579055
$ws.Cells.Item(4,13).Value = 0.1408326666666667
$ws.Cells.Item(4,14).Value = 0.422498
$ws.Cells.Item(4,15).Value = 0.0430631597358264
$ws.Cells.Item(4,16).Value = 0.04306315973582638
$ws.Cells.Item(4,17).Value = 0.2249864755257778
$ws.Cells.Item(4,18).Value = 2.024878279732
$ws.Cells.Item(4,19).Value = 0.0009884205340005912
$ws.Cells.Item(4,20).Value = 0.0009884205340005908
$ws.Cells.Item(5,7).Value = 1.597544666666667
$ws.Cells.Item(5,8).Value = 4.792634
$ws.Cells.Item(5,9).Value = 0.02295281024579055
$ws.Cells.Item(5,10).Value = 0.02295281024579055
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.072925
$ws.Cells.Item(5,14).Value = 0.218775
$ws.Cells.Item(5,15).Value = 0.02229866832791023
$ws.Cells.Item(5,16).Value = 0.02229866832791023
$ws.Cells.Item(5,17).Value = 0.1165009448166667
$ws.Cells.Item(5,18).Value = 1.04850850335
$ws.Cells.Item(5,19).Value = 0.0005118171028643433
$ws.Cells.Item(5,20).Value = 0.0005118171028643432
$ws.Cells.Item(6,9).Value = 0.7618977001376412
$ws.Cells.Item(6,10).Value = 0.7618977001376412
$ws.Cells.Item(6,13).Value = 2.680128333333334
$ws.Cells.Item(6,14).Value = 8.040385000000001
$ws.Cells.Item(6,15).Value = 0.8195172133182701
$ws.Cells.Item(6,16).Value = 0.81951721331827
$ws.Cells.Item(6,17).Value = 142.1246459941083
$ws.Cells.Item(6,18).Value = 1279.121813946975
$ws.Cells.Item(6,19).Value = 0.6243882800503987
$ws.Cells.Item(6,20).Value = 0.6243882800503986
$ws.Cells.Item(7,9).Value = 0.7618977001376412
$ws.Cells.Item(7,10).Value = 0.7618977001376412
$ws.Cells.Item(7,13).Value = 0.3764886666666666
$ws.Cells.Item(7,15).Value = 0.1151209586179932
$ws.Cells.Item(7,16).Value = 0.1151209586179932
$ws.Cells.Item(7,17).Value = 19.96483444665666
$ws.Cells.Item(7,19).Value = 0.08771039360868957
$ws.Cells.Item(7,20).Value = 0.08771039360868955
$ws.Cells.Item(8,9).Value = 0.7618977001376412
$ws.Cells.Item(8,10).Value = 0.7618977001376412
$ws.Cells.Item(8,13).Value = 0.1408326666666667
$ws.Cells.Item(8,14).Value = 0.422498
$ws.Cells.Item(8,15).Value = 0.0430631597358264
$ws.Cells.Item(8,16).Value = 0.04306315973582638
$ws.Cells.Item(8,17).Value = 7.468221818136668
$ws.Cells.Item(8,18).Value = 67.21399636323
$ws.Cells.Item(8,19).Value = 0.03280972236338601
$ws.Cells.Item(8,20).Value = 0.03280972236338599
$ws.Cells.Item(9,9).Value = 0.7618977001376412
$ws.Cells.Item(9,10).Value = 0.7618977001376412
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.072925
$ws.Cells.Item(9,14).Value = 0.218775
$ws.Cells.Item(9,15).Value = 0.02229866832791023
$ws.Cells.Item(9,16).Value = 0.02229866832791023
$ws.Cells.Item(9,17).Value = 3.867143106625
$ws.Cells.Item(9,18).Value = 34.804287959625
$ws.Cells.Item(9,19).Value = 0.01698930411516687
$ws.Cells.Item(9,20).Value = 0.01698930411516687
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.1019876666666667
$ws.Cells.Item(10,8).Value = 0.305963
$ws.Cells.Item(10,9).Value = 0.001465313370733674
$ws.Cells.Item(10,10).Value = 0.001465313370733675
$ws.Cells.Item(10,13).Value = 2.680128333333334
$ws.Cells.Item(10,14).Value = 8.040385000000001
$ws.Cells.Item(10,15).Value = 0.8195172133182701
$ws.Cells.Item(10,16).Value = 0.81951721331827
$ws.Cells.Item(10,17).Value = 0.2733400350838889
$ws.Cells.Item(10,18).Value = 2.460060315755
$ws.Cells.Item(10,19).Value = 0.001200849530221662
$ws.Cells.Item(10,20).Value = 0.001200849530221662
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.1019876666666667
$ws.Cells.Item(11,8).Value = 0.305963
$ws.Cells.Item(11,9).Value = 0.001465313370733674
$ws.Cells.Item(11,10).Value = 0.001465313370733675
$ws.Cells.Item(11,13).Value = 0.3764886666666666
$ws.Cells.Item(11,15).Value = 0.1151209586179932
$ws.Cells.Item(11,16).Value = 0.1151209586179932
$ws.Cells.Item(11,17).Value = 0.03839720063977777
$ws.Cells.Item(11,18).Value = 0.3455748057579999
$ws.Cells.Item(11,19).Value = 0.0001686882799146234
$ws.Cells.Item(11,20).Value = 0.0001686882799146234
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.1019876666666667
$ws.Cells.Item(12,8).Value = 0.305963
$ws.Cells.Item(12,9).Value = 0.001465313370733674
$ws.Cells.Item(12,10).Value = 0.001465313370733675
$ws.Cells.Item(12,13).Value = 0.1408326666666667
$ws.Cells.Item(12,14).Value = 0.422498
$ws.Cells.Item(12,15).Value = 0.0430631597358264
$ws.Cells.Item(12,16).Value = 0.04306315973582638
$ws.Cells.Item(12,17).Value = 0.01436319506377778
$ws.Cells.Item(12,18).Value = 0.129268755574
$ws.Cells.Item(12,19).Value = [double]"6.310102374694643E-05"
$ws.Cells.Item(12,20).Value = [double]"6.310102374694642E-05"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.1019876666666667
$ws.Cells.Item(13,8).Value = 0.305963
$ws.Cells.Item(13,9).Value = 0.001465313370733674
$ws.Cells.Item(13,10).Value = 0.001465313370733675
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.072925
$ws.Cells.Item(13,14).Value = 0.218775
$ws.Cells.Item(13,15).Value = 0.02229866832791023
$ws.Cells.Item(13,16).Value = 0.02229866832791023
$ws.Cells.Item(13,17).Value = 0.007437450591666666
$ws.Cells.Item(13,18).Value = 0.06693705532499999
$ws.Cells.Item(13,19).Value = [double]"3.267453685044238E-05"
$ws.Cells.Item(13,20).Value = [double]"3.267453685044238E-05"
$ws.Cells.Item(14,7).Value = 14.16318533333333
$ws.Cells.Item(14,8).Value = 42.489556
$ws.Cells.Item(14,9).Value = 0.2034903387773594
$ws.Cells.Item(14,10).Value = 0.2034903387773594
$ws.Cells.Item(14,13).Value = 2.680128333333334
$ws.Cells.Item(14,14).Value = 8.040385000000001
$ws.Cells.Item(14,15).Value = 0.8195172133182701
$ws.Cells.Item(14,16).Value = 0.81951721331827
$ws.Cells.Item(14,17).Value = 37.95915430211778
$ws.Cells.Item(14,18).Value = 341.63238871906
$ws.Cells.Item(14,19).Value = 0.1667638353720123
$ws.Cells.Item(14,20).Value = 0.1667638353720123
$ws.Cells.Item(15,7).Value = 14.16318533333333
$ws.Cells.Item(15,8).Value = 42.489556
$ws.Cells.Item(15,9).Value = 0.2034903387773594
$ws.Cells.Item(15,10).Value = 0.2034903387773594
$ws.Cells.Item(15,13).Value = 0.3764886666666666
$ws.Cells.Item(15,15).Value = 0.1151209586179932
$ws.Cells.Item(15,16).Value = 0.1151209586179932
$ws.Cells.Item(15,17).Value = 5.332278761899555
$ws.Cells.Item(15,18).Value = 47.990508857096
$ws.Cells.Item(15,19).Value = 0.02342600286954981
$ws.Cells.Item(15,20).Value = 0.02342600286954981
$ws.Cells.Item(16,7).Value = 14.16318533333333
$ws.Cells.Item(16,8).Value = 42.489556
$ws.Cells.Item(16,9).Value = 0.2034903387773594
$ws.Cells.Item(16,10).Value = 0.2034903387773594
$ws.Cells.Item(16,13).Value = 0.1408326666666667
$ws.Cells.Item(16,14).Value = 0.422498
$ws.Cells.Item(16,15).Value = 0.0430631597358264
$ws.Cells.Item(16,16).Value = 0.04306315973582638
$ws.Cells.Item(16,17).Value = 1.994639158987556
$ws.Cells.Item(16,18).Value = 17.951752430888
$ws.Cells.Item(16,19).Value = 0.008762936963466857
$ws.Cells.Item(16,20).Value = 0.008762936963466856
$ws.Cells.Item(17,7).Value = 14.16318533333333
$ws.Cells.Item(17,8).Value = 42.489556
$ws.Cells.Item(17,9).Value = 0.2034903387773594
$ws.Cells.Item(17,10).Value = 0.2034903387773594
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.072925
$ws.Cells.Item(17,14).Value = 0.218775
$ws.Cells.Item(17,15).Value = 0.02229866832791023
$ws.Cells.Item(17,16).Value = 0.02229866832791023
$ws.Cells.Item(17,17).Value = 1.032850290433333
$ws.Cells.Item(17,18).Value = 9.2956526139
$ws.Cells.Item(17,19).Value = 0.004537563572330428
$ws.Cells.Item(17,20).Value = 0.004537563572330428
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 0.709504
$ws.Cells.Item(18,8).Value = 2.128512
$ws.Cells.Item(18,9).Value = 0.01019383746847519
$ws.Cells.Item(18,10).Value = 0.01019383746847519
$ws.Cells.Item(18,13).Value = 2.680128333333334
$ws.Cells.Item(18,14).Value = 8.040385000000001
$ws.Cells.Item(18,15).Value = 0.8195172133182701
$ws.Cells.Item(18,16).Value = 0.81951721331827
$ws.Cells.Item(18,17).Value = 1.901561773013334
$ws.Cells.Item(18,18).Value = 17.11405595712
$ws.Cells.Item(18,19).Value = 0.008354025275184158
$ws.Cells.Item(18,20).Value = 0.008354025275184158
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 0.709504
$ws.Cells.Item(19,8).Value = 2.128512
$ws.Cells.Item(19,9).Value = 0.01019383746847519
$ws.Cells.Item(19,10).Value = 0.01019383746847519
$ws.Cells.Item(19,13).Value = 0.3764886666666666
$ws.Cells.Item(19,15).Value = 0.1151209586179932
$ws.Cells.Item(19,16).Value = 0.1151209586179932
$ws.Cells.Item(19,17).Value = 0.2671202149546666
$ws.Cells.Item(19,18).Value = 2.404081934592
$ws.Cells.Item(19,19).Value = 0.001173524341366881
$ws.Cells.Item(19,20).Value = 0.001173524341366881
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 0.709504
$ws.Cells.Item(20,8).Value = 2.128512
$ws.Cells.Item(20,9).Value = 0.01019383746847519
$ws.Cells.Item(20,10).Value = 0.01019383746847519
$ws.Cells.Item(20,13).Value = 0.1408326666666667
$ws.Cells.Item(20,14).Value = 0.422498
$ws.Cells.Item(20,15).Value = 0.0430631597358264
$ws.Cells.Item(20,16).Value = 0.04306315973582638
$ws.Cells.Item(20,17).Value = 0.09992134033066669
$ws.Cells.Item(20,18).Value = 0.8992920629760002
$ws.Cells.Item(20,19).Value = 0.0004389788512259994
$ws.Cells.Item(20,20).Value = 0.0004389788512259993
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 0.709504
$ws.Cells.Item(21,8).Value = 2.128512
$ws.Cells.Item(21,9).Value = 0.01019383746847519
$ws.Cells.Item(21,10).Value = 0.01019383746847519
$ws.Cells.Item(21,11).Value = 2
$ws.Cells.Item(21,12).Value = 0.6666666666666666
$ws.Cells.Item(21,13).Value = 0.072925
$ws.Cells.Item(21,14).Value = 0.218775
$ws.Cells.Item(21,15).Value = 0.02229866832791023
$ws.Cells.Item(21,16).Value = 0.02229866832791023
$ws.Cells.Item(21,17).Value = 0.0517405792
$ws.Cells.Item(21,18).Value = 0.4656652128000001
$ws.Cells.Item(21,19).Value = 0.0002273090006981524
$ws.Cells.Item(21,20).Value = 0.0002273090006981524
